# LOM3240.xlsx update
#  - "Ativacao:" date bumped from 01/01/2012 -> 01/01/2023 (row 8 / row 15, B & C columns)
#  - English "Objectives:" (row 11) gets its B/C text filled in
#  - English "Short syllabus:" (row 14) gets its B/C text filled in
#  - English "Syllabus:" (row 16) gets its B/C text filled in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Ativacao / Activation date: 01/01/2012 -> 01/01/2023
#    Kept as literal text (the workbook stores it as a shared string, not a
#    date serial), so the value is entered with a leading apostrophe to stop
#    it being auto-parsed into a date by the COM layer.
# ---------------------------------------------------------------------------
$newDate = "'01/01/2023"

$ws.Range("B8").Value = $newDate
$ws.Range("C8").Value = $newDate
$ws.Range("B15").Value = $newDate
$ws.Range("C15").Value = $newDate

# ---------------------------------------------------------------------------
# 2) Objectives: (row 11) - add the English translation in columns B and C
# ---------------------------------------------------------------------------
$objectivesEn = "To present an overview of the chemistry of elements and their compounds, emphasizing the correlations between physical and chemical properties with structural and binding aspects, methods of obtaining them in laboratory and industry, in addition to the main properties and applications."

$ws.Range("B11").Font.Bold = $false
$ws.Range("B11").WrapText = $true
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").Value = $objectivesEn
$ws.Range("C11").Value = $objectivesEn

# ---------------------------------------------------------------------------
# 3) Short syllabus: (row 14) - add the English translation in columns B and C
# ---------------------------------------------------------------------------
$shortSyllabusEn = "Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes."

$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").Value = $shortSyllabusEn
$ws.Range("C14").Value = $shortSyllabusEn

# ---------------------------------------------------------------------------
# 4) Syllabus: (row 16) - add the English translation in columns B and C
# ---------------------------------------------------------------------------
$syllabusEn = "Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials."

$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("B16").Value = $syllabusEn
$ws.Range("C16").Value = $syllabusEn

Write-Output "Done."
